# Apply the parameter updates for the "base" and "v4" parameterization
# comparison described in the commit message:
#   "add parameter comparison code, change base and v4 parameterizations"
#
# Concretely this means updating a handful of values in the stakeholder /
# parameter matrix on Sheet1, and leaving the user's selection on the last
# edited cell (J28), matching the state the workbook was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Drinking Water Division (SWRCB)" row - base parameterization value
# changed from 2 to 1.
$ws.Range("B17").Value = 1

# "Local Water Boards" row - v4 parameterization value changed from 2 to 1.
$ws.Range("B22").Value = 1

# "CV SALTS management zones" row - fill in the comparison flags for the
# remaining stakeholder columns (C:F) to -1, matching columns G:H which
# were already set.
$ws.Range("C23:F23").Value = -1

# Leave the selection where the author left it when they saved the file.
[void]$ws.Range("J28").Select()
